$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency market data (price & 1h volume change)
# Force text number format per-cell so values are stored as text (matching original inlineStr cells)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "242.82"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.01%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "29.70"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "11.98%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.131"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.24%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05653"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.25%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.493"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.22%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8273"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.27%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8620"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.69%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1328"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.19%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06920"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.12%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.02858"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.96%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09388"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.08%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001524"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.38%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04150"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-9.48%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006009"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-94.00%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006168"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.67%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.520"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-3.15%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.79%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.218"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.58%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.03249"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.95%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.630"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.86%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.04%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-3.00%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004445"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-1.49%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001179"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "22.86%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "0.57%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03704"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.75%"
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005746"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-6.70%"
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1053"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.21%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002310"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.73%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009657"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "9.12%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005109"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.39%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.04%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1010"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-7.34%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002582"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "1.18%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.04%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001999"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.04%"
